$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet: reset capital / active-strategy counters to 0
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0    # Initial Capital
$summary.Range("B3").Value = 0    # Current Capital
$summary.Range("B11").Value = 0   # Active Strategies

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet: drop every strategy row, keep only the header
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("A2:G16").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 3. All Trades sheet: clear out the now-stale MarketMaking rows (5 & 6) and
#    append the new TestStrategy trade (7)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Trade #5 (row 6) - exit price now recorded as 0, extra fields cleared
$allTrades.Range("G6").Value = 0
$allTrades.Range("K6").Value = ""
$allTrades.Range("L6").Value = ""
$allTrades.Range("M6").Value = ""
$allTrades.Range("N6").Value = ""
$allTrades.Range("O6").Value = ""
$allTrades.Range("Q6").Value = ""

# Trade #6 (row 7) - exit price now recorded as 0, extra fields cleared
$allTrades.Range("G7").Value = 0
$allTrades.Range("K7").Value = ""
$allTrades.Range("L7").Value = ""
$allTrades.Range("M7").Value = ""
$allTrades.Range("N7").Value = ""
$allTrades.Range("O7").Value = ""
$allTrades.Range("Q7").Value = ""

# Trade #7 (new row 8) - TestStrategy entry
$allTrades.Range("A8").Value = 7
$allTrades.Range("B8").NumberFormat = "@"
$allTrades.Range("B8").Value = "2026-02-18"
$allTrades.Range("C8").Value = "10:54:40"
$allTrades.Range("D8").Value = "TestStrategy"
$allTrades.Range("E8").Value = "UP"
$allTrades.Range("F8").Value = 0.5
$allTrades.Range("H8").Value = "OPEN"
$allTrades.Range("I8").Value = 0
$allTrades.Range("J8").Value = 0
$allTrades.Range("K8").Value = 100
$allTrades.Range("L8").Value = 0
$allTrades.Range("M8").Value = 0
$allTrades.Range("N8").Value = 0.8
$allTrades.Range("O8").Value = "Test entry"
$allTrades.Range("Q8").Value = 0

# ---------------------------------------------------------------------------
# 4. Per-strategy sheet: "MarketMaking" becomes "TestStrategy".
#    Its old second trade (row 3) is dropped, and the remaining row is
#    updated in place to reflect the new TestStrategy trade #7.
# ---------------------------------------------------------------------------
$strategySheet = $wb.Worksheets.Item("MarketMaking")
$strategySheet.Range("A3:Q3").EntireRow.Delete()

$strategySheet.Range("A2").Value = 7
$strategySheet.Range("C2").Value = "10:54:40"
$strategySheet.Range("D2").Value = "TestStrategy"
$strategySheet.Range("E2").Value = "UP"
$strategySheet.Range("F2").Value = 0.5
$strategySheet.Range("N2").Value = 0.8
$strategySheet.Range("O2").Value = "Test entry"

$strategySheet.Name = "TestStrategy"
